$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric survey data (columns A-F) and the category label (column H)
# for rows 1-6 to reflect the 2021 questionnaire results.

$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 4
$ws.Range("D1").Value = 12
$ws.Range("E1").Value = 1
$ws.Range("F1").Value = 0
$ws.Range("H1").Value = "D"

$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("H2").Value = "B-C"

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("H3").Value = "B"

$ws.Range("A4").Value = 9
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("H4").Value = "A-B"

$ws.Range("A5").Value = 11
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("H5").Value = "A"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0
$ws.Range("H6").Value = "A-E"

# Row 7 no longer exists in the 2021 data set - delete it entirely.
$ws.Rows.Item(7).Delete()

# Restore the previously-active cell selection.
$ws.Range("H7").Select()
